$d = $word.ActiveDocument

# --- 1. Title paragraph ("Feb. 18-21 To-dos"): drop the explicit 24-half-point
#        size, keeping Bold + center alignment intact. ClearFormatting wipes
#        run AND paragraph-mark direct formatting, so re-apply Bold after.
$titlePar = $d.Paragraphs.Item(1)
$titlePar.Range.Select()
$word.Selection.ClearFormatting()
$word.Selection.Font.Bold = $true
$titlePar.Alignment = 1

# --- 2. Insert the two new "camera" to-dos (italic, matching the surrounding
#        list style) before doing the size sweep, so they naturally pick up
#        sz=20 like every other to-do line.

# "-camera mounting" goes right after "-pot", before the "Bumpers" heading.
$potPar = $d.Paragraphs.Item(11)
$potPar.Range.InsertParagraphAfter()
$cameraMountPar = $d.Paragraphs.Item(12)
$cameraMountPar.Range.Text = "-camera mounting"
$cameraMountPar.Range.Font.Italic = $true

# "-camera functionality" goes right after "-make the elbow go", before the
# "Misc." heading. (Index shifted by +1 because of the insert above.)
$elbowPar = $d.Paragraphs.Item(20)
$elbowPar.Range.InsertParagraphAfter()
$cameraFuncPar = $d.Paragraphs.Item(21)
$cameraFuncPar.Range.Text = "-camera functionality"
$cameraFuncPar.Range.Font.Italic = $true

# --- 3. Every other paragraph (2 through the new end of the to-do list, i.e.
#        everything except the title and the trailing blank paragraph) gets
#        a 10pt (sz=20) font so the to-do body shrinks to fit more lines.
$total = $d.Paragraphs.Count
for ($i = 2; $i -le ($total - 1); $i++) {
    $d.Paragraphs.Item($i).Range.Font.Size = 10
}
